$wb = $excel.ActiveWorkbook
$count = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($count)
$ws = $wb.Worksheets.Add([System.Type]::Missing, $lastSheet)
$ws.Name = "2025-08-16"

# Copy header row (values + formatting) from the previous day sheet
$srcHeader = $wb.Worksheets.Item("2025-08-15").Range("A1:D1")
$srcHeader.Copy($ws.Range("A1:D1"))

# Populate ranking data rows 2-51
$ws.Cells.Item(2, 1).Value = 1
$ws.Cells.Item(2, 2).Value = '宇崎ちゃんは遊びたい！'
$ws.Cells.Item(2, 3).Value = '丈(著者)'
$ws.Cells.Item(2, 4).Value = '第126話'
$ws.Cells.Item(3, 1).Value = 2
$ws.Cells.Item(3, 2).Value = 'よくわからないけれど異世界に転生していたようです'
$ws.Cells.Item(3, 3).Value = '内々けやき あし カオミン'
$ws.Cells.Item(3, 4).Value = '第137話 よくわからないけれど脱出するみたいです（２）'
$ws.Cells.Item(4, 1).Value = 3
$ws.Cells.Item(4, 2).Value = '小林さんちのメイドラゴン'
$ws.Cells.Item(4, 3).Value = 'クール教信者'
$ws.Cells.Item(4, 4).Value = '第148話'
$ws.Cells.Item(5, 1).Value = 4
$ws.Cells.Item(5, 2).Value = 'くらいあの子としたいこと'
$ws.Cells.Item(5, 3).Value = '碇マナツ(著者)'
$ws.Cells.Item(5, 4).Value = '第82話'
$ws.Cells.Item(6, 1).Value = 5
$ws.Cells.Item(6, 2).Value = '治癒魔法の間違った使い方 ~戦場を駆ける回復要員~'
$ws.Cells.Item(6, 3).Value = '九我山レキ(漫画) くろかた(原作) ＫｅＧ(キャラクター原案)'
$ws.Cells.Item(6, 4).Value = '第81話その3'
$ws.Cells.Item(7, 1).Value = 6
$ws.Cells.Item(7, 2).Value = '理想のヒモ生活'
$ws.Cells.Item(7, 3).Value = '日月ネコ(漫画) 渡辺恒彦（ヒーロー文庫／イマジカインフォス）(原作) 文倉十(キャラクター原案)'
$ws.Cells.Item(7, 4).Value = '第86話　その3'
$ws.Cells.Item(8, 1).Value = 7
$ws.Cells.Item(8, 2).Value = '追放されたチート付与魔術師は 気ままなセカンドライフを謳歌する。'
$ws.Cells.Item(8, 3).Value = '六志麻あさ 業務用餅 kisui'
$ws.Cells.Item(8, 4).Value = '第７０話'
$ws.Cells.Item(9, 1).Value = 8
$ws.Cells.Item(9, 2).Value = '最強勇者パーティーは愛が知りたい'
$ws.Cells.Item(9, 3).Value = '山田肌襦袢'
$ws.Cells.Item(9, 4).Value = '第29話「きみがきみであればいい」'
$ws.Cells.Item(10, 1).Value = 9
$ws.Cells.Item(10, 2).Value = 'えろいことするために巨乳美少女奴隷を買ったはずが、お師匠さまと慕われて思った通りにいかなくなる話'
$ws.Cells.Item(10, 3).Value = '佐藤36(作画) 煮豆シューター(原作)'
$ws.Cells.Item(10, 4).Value = '第2話後半'
$ws.Cells.Item(11, 1).Value = 10
$ws.Cells.Item(11, 2).Value = '序盤で死ぬ最強のサブキャラに転生したので、ゲーム知識で無双する'
$ws.Cells.Item(11, 3).Value = '作画：マエD 原作：新人'
$ws.Cells.Item(11, 4).Value = '第5話(4)'
$ws.Cells.Item(12, 1).Value = 11
$ws.Cells.Item(12, 2).Value = '最強の少年聖騎士、転生者を狩る'
$ws.Cells.Item(12, 3).Value = '作画：御塩 原作：宇奈木ユラ'
$ws.Cells.Item(12, 4).Value = '第7話(1)'
$ws.Cells.Item(13, 1).Value = 12
$ws.Cells.Item(13, 2).Value = 'みつばものがたり 呪いの少女と死の輪舞《ロンド》'
$ws.Cells.Item(13, 3).Value = '堤りん(漫画) 七沢またり(原作) EURA(キャラクター原案)'
$ws.Cells.Item(13, 4).Value = '第11話：勝利の美酒'
$ws.Cells.Item(14, 1).Value = 13
$ws.Cells.Item(14, 2).Value = 'ゲーセン少女と異文化交流'
$ws.Cells.Item(14, 3).Value = '安原宏和(著者)'
$ws.Cells.Item(14, 4).Value = '第128話'
$ws.Cells.Item(15, 1).Value = 14
$ws.Cells.Item(15, 2).Value = '不徳のギルド'
$ws.Cells.Item(15, 3).Value = '河添太一'
$ws.Cells.Item(15, 4).Value = '第９７話：立派に育った所'
$ws.Cells.Item(16, 1).Value = 15
$ws.Cells.Item(16, 2).Value = 'ヤンデレかと思ったらもっとヤベー女だった'
$ws.Cells.Item(16, 3).Value = '八木戸マト'
$ws.Cells.Item(16, 4).Value = '第71話　奪い返しにきたヤンデレ彼女'
$ws.Cells.Item(17, 1).Value = 16
$ws.Cells.Item(17, 2).Value = '王都の外れの錬金術師 ～ハズレ職業だったので、のんびりお店経営します～'
$ws.Cells.Item(17, 3).Value = 'あさなや(著者) yocco(原作) 純粋(キャラクター原案)'
$ws.Cells.Item(17, 4).Value = 'element.50'
$ws.Cells.Item(18, 1).Value = 17
$ws.Cells.Item(18, 2).Value = '路地裏で拾った女の子がバッドエンド後の乙女ゲームのヒロインだった件'
$ws.Cells.Item(18, 3).Value = 'カボチャマスク(原作) 樋乃えなが(作画) へいろー(キャラクター原案)'
$ws.Cells.Item(18, 4).Value = '第1話'
$ws.Cells.Item(19, 1).Value = 18
$ws.Cells.Item(19, 2).Value = '修羅幼女の英雄譚～半端者と言われた傭兵、幼女に転生して成り上がる～'
$ws.Cells.Item(19, 3).Value = '作画：むらたん 原作：沙城流'
$ws.Cells.Item(19, 4).Value = '第8話(3)'
$ws.Cells.Item(20, 1).Value = 19
$ws.Cells.Item(20, 2).Value = 'お前妹じゃなくて許嫁だったのかよ!?'
$ws.Cells.Item(20, 3).Value = '湯猫子(漫画) 未来人A(原作)'
$ws.Cells.Item(20, 4).Value = '第29話'
$ws.Cells.Item(21, 1).Value = 20
$ws.Cells.Item(21, 2).Value = '最強で最速の無限レベルアップ ～スキル【経験値1000倍】と【レベルフリー】でレベル上限の枷が外れた俺は無双する～'
$ws.Cells.Item(21, 3).Value = 'シオヤマ琴 鳥羽田 航 トモゼロ'
$ws.Cells.Item(21, 4).Value = '休載イラスト'
$ws.Cells.Item(22, 1).Value = 21
$ws.Cells.Item(22, 2).Value = '役目を果たした日陰の勇者は、辺境で自由に生きていきます'
$ws.Cells.Item(22, 3).Value = '船野真帆(作画) 丘野優(原作) 布施龍太(キャラクター原案)'
$ws.Cells.Item(22, 4).Value = '第5話前半'
$ws.Cells.Item(23, 1).Value = 22
$ws.Cells.Item(23, 2).Value = 'なぜかS級美女達の話題に俺があがる件'
$ws.Cells.Item(23, 3).Value = 'ジョN(著者) 脇岡こなつ(原作) magako(キャラクター原案)'
$ws.Cells.Item(23, 4).Value = '最終話-1'
$ws.Cells.Item(24, 1).Value = 23
$ws.Cells.Item(24, 2).Value = '農学博士の異世界無双～禁忌の知識で築くモンスター娘ハーレム～'
$ws.Cells.Item(24, 3).Value = 'インド僧(原作) ヤスウミ(作画)'
$ws.Cells.Item(24, 4).Value = '第25話'
$ws.Cells.Item(25, 1).Value = 24
$ws.Cells.Item(25, 2).Value = '異世界はスマートフォンとともに。'
$ws.Cells.Item(25, 3).Value = 'そと(漫画) 冬原パトラ(原作) 兎塚エイジ(キャラクター原案)'
$ws.Cells.Item(25, 4).Value = 'EPISODE:103‐②'
$ws.Cells.Item(26, 1).Value = 25
$ws.Cells.Item(26, 2).Value = '顔に出ない柏田さんと顔に出る太田君＋'
$ws.Cells.Item(26, 3).Value = '東ふゆ(著者)'
$ws.Cells.Item(26, 4).Value = '第32話 田淵さんと愛田さんのバトル'
$ws.Cells.Item(27, 1).Value = 26
$ws.Cells.Item(27, 2).Value = '底辺冒険者だけど魔法を極めてみることにした ～無能スキルから神スキルに進化した【魔法創造】と【アイテム作成】で無双する～'
$ws.Cells.Item(27, 3).Value = '蒼乃白兎 坂野杏梨 かわく'
$ws.Cells.Item(27, 4).Value = '第45話(前編) 反逆者'
$ws.Cells.Item(28, 1).Value = 27
$ws.Cells.Item(28, 2).Value = 'ダウナーお姉さんは遊びたい'
$ws.Cells.Item(28, 3).Value = '山鷹景'
$ws.Cells.Item(28, 4).Value = '第18話'
$ws.Cells.Item(29, 1).Value = 28
$ws.Cells.Item(29, 2).Value = 'ノロマ魔法と呼ばれた魔法使いは重力魔法で無双する　～まだ重力の概念のない世界にて、少年は万有引力の王となる～'
$ws.Cells.Item(29, 3).Value = '神原絵理華(漫画) 一森一輝(原作)'
$ws.Cells.Item(29, 4).Value = '描き下ろしイラスト公開！⑫'
$ws.Cells.Item(30, 1).Value = 29
$ws.Cells.Item(30, 2).Value = 'ゴミ以下だと追放された使用人、実は前世賢者です　～史上最強の賢者、世界最高峰の学園に通う～'
$ws.Cells.Item(30, 3).Value = '原作：夜分長文 漫画：矢部利恩 キャラクター原案：蔓木鋼音'
$ws.Cells.Item(30, 4).Value = '第15話 因縁の再会（１）'
$ws.Cells.Item(31, 1).Value = 30
$ws.Cells.Item(31, 2).Value = 'パワハラギルマスをぶん殴ってブラック聖剣ギルドをクビになったので、辺境で聖剣工房を開くことにした'
$ws.Cells.Item(31, 3).Value = 'だいたいねむい(原作) まお(漫画)'
$ws.Cells.Item(31, 4).Value = '第10話②'
$ws.Cells.Item(32, 1).Value = 31
$ws.Cells.Item(32, 2).Value = 'スキルがなければレベルを上げる～９９がカンストの世界でレベル800万からスタート～'
$ws.Cells.Item(32, 3).Value = '倉橋ユウス(漫画) 岡沢六十四(原作)'
$ws.Cells.Item(32, 4).Value = '第52話③'
$ws.Cells.Item(33, 1).Value = 32
$ws.Cells.Item(33, 2).Value = '辺境モブ貴族のウチに嫁いできた悪役令嬢が、めちゃくちゃできる良い嫁なんだが？'
$ws.Cells.Item(33, 3).Value = 'tera(原作) 朝倉はやて(作画) 徹田(キャラクター原案)'
$ws.Cells.Item(33, 4).Value = '第10話'
$ws.Cells.Item(34, 1).Value = 33
$ws.Cells.Item(34, 2).Value = '傭兵団の愛し子 ～死にかけ孤児は最強師匠たちに育てられる～'
$ws.Cells.Item(34, 3).Value = '柿野レイ(漫画) 天野雪人(原作) 黒井ススム(キャラクター原案)'
$ws.Cells.Item(34, 4).Value = '第7話後編：問題児シャルノア'
$ws.Cells.Item(35, 1).Value = 34
$ws.Cells.Item(35, 2).Value = '限界超えの天賦は、転生者にしか扱えない ― オーバーリミット・スキルホルダー ―'
$ws.Cells.Item(35, 3).Value = '長月みそか(漫画) 三上康明(原作) 大槍葦人(キャラクターデザイン)'
$ws.Cells.Item(35, 4).Value = '第3章［14］後半'
$ws.Cells.Item(36, 1).Value = 35
$ws.Cells.Item(36, 2).Value = '俺の『全自動支援（フルオートバフ）』で仲間たちが世界最強 ～そこにいるだけ無自覚無双～'
$ws.Cells.Item(36, 3).Value = 'IプルT(著者) epina(原作) 片倉響(キャラクター原案)'
$ws.Cells.Item(36, 4).Value = '第３４話「夢、再び」'
$ws.Cells.Item(37, 1).Value = 36
$ws.Cells.Item(37, 2).Value = '勇者パーティから追い出された不遇職【罠士】、ユニークスキル【矢印】で最強になる'
$ws.Cells.Item(37, 3).Value = '作画：たつひこ 原作：白石 有希'
$ws.Cells.Item(37, 4).Value = '第8話(3)'
$ws.Cells.Item(38, 1).Value = 37
$ws.Cells.Item(38, 2).Value = '道にスライムが捨てられていたから連れて帰りました ～おじさんとスライムのほのぼの冒険ライフ～'
$ws.Cells.Item(38, 3).Value = 'めぐお(漫画) イコ(原作) いもいち(キャラクター原案)'
$ws.Cells.Item(38, 4).Value = '第2話-2'
$ws.Cells.Item(39, 1).Value = 38
$ws.Cells.Item(39, 2).Value = '俺の愛娘は悪役令嬢'
$ws.Cells.Item(39, 3).Value = 'かわもり かぐら(原作) ほづみりや(漫画) 縞(キャラクター原案)'
$ws.Cells.Item(39, 4).Value = '第5話-1'
$ws.Cells.Item(40, 1).Value = 39
$ws.Cells.Item(40, 2).Value = '義妹生活'
$ws.Cells.Item(40, 3).Value = '三河ごーすと(原作) 奏ユミカ(漫画) Hiten(キャラクター原案)'
$ws.Cells.Item(40, 4).Value = '第31話-2'
$ws.Cells.Item(41, 1).Value = 40
$ws.Cells.Item(41, 2).Value = 'ブレイド＆バスタード'
$ws.Cells.Item(41, 3).Value = '漫画/楓月 誠 原作/蝸牛くも キャラクター原案/so-bin'
$ws.Cells.Item(41, 4).Value = '第11話（2）'
$ws.Cells.Item(42, 1).Value = 41
$ws.Cells.Item(42, 2).Value = 'やっぱ人間やめて正解だわ'
$ws.Cells.Item(42, 3).Value = '偽BEなんとか'
$ws.Cells.Item(42, 4).Value = '第16話（後編）　ちんちんな二人'
$ws.Cells.Item(43, 1).Value = 42
$ws.Cells.Item(43, 2).Value = 'アイツノカノジョ'
$ws.Cells.Item(43, 3).Value = '肉丸'
$ws.Cells.Item(43, 4).Value = '第55話'
$ws.Cells.Item(44, 1).Value = 43
$ws.Cells.Item(44, 2).Value = 'うちの清楚系委員長がかつて中二病アイドルだったことを俺だけが知っている。'
$ws.Cells.Item(44, 3).Value = '三上こた こばやし少女 寝子空兄 ゆがー'
$ws.Cells.Item(44, 4).Value = '第1話　災禍の悪夢'
$ws.Cells.Item(45, 1).Value = 44
$ws.Cells.Item(45, 2).Value = 'ぽんドロイド！ はまさん'
$ws.Cells.Item(45, 3).Value = 'はれやまはれぞう(著者)'
$ws.Cells.Item(45, 4).Value = '第6話'
$ws.Cells.Item(46, 1).Value = 45
$ws.Cells.Item(46, 2).Value = 'バーサス'
$ws.Cells.Item(46, 3).Value = '原作：ONE 漫画：あずま京太郎 構成：bose'
$ws.Cells.Item(46, 4).Value = '第27話 幸せの在り処（1）'
$ws.Cells.Item(47, 1).Value = 46
$ws.Cells.Item(47, 2).Value = '落ちこぼれだった兄が実は最強 ～史上最強の勇者は転生し、学園で無自覚に無双する～'
$ws.Cells.Item(47, 3).Value = '村上よしゆき 茨木野 あるてら'
$ws.Cells.Item(47, 4).Value = '第４１話　勇者、人魚王国を救い、歓迎される。あと、六邪神将が、全員来る（３）'
$ws.Cells.Item(48, 1).Value = 47
$ws.Cells.Item(48, 2).Value = '隠居暮らしのおっさん、女王陛下の剣となる'
$ws.Cells.Item(48, 3).Value = '漫画/半二合 原作/天酒之瓢 キャラクター原案/みことあけみ'
$ws.Cells.Item(48, 4).Value = '第6話（1）'
$ws.Cells.Item(49, 1).Value = 48
$ws.Cells.Item(49, 2).Value = 'ギルドを追放された回復術士、実は魔力無限だったので規格外の回復魔法で伝説となる'
$ws.Cells.Item(49, 3).Value = '漫画：坂下コウ 原作：霞杏檎'
$ws.Cells.Item(49, 4).Value = '第5話(2)'
$ws.Cells.Item(50, 1).Value = 49
$ws.Cells.Item(50, 2).Value = '異世界のんびり農家の日常'
$ws.Cells.Item(50, 3).Value = 'ユウズィ(著者) 内藤騎之介(原作) やすも(キャラクター原案)'
$ws.Cells.Item(50, 4).Value = '収穫その37'
$ws.Cells.Item(51, 1).Value = 50
$ws.Cells.Item(51, 2).Value = 'やめてくれ、強いのは俺じゃなくて剣なんだ……！'
$ws.Cells.Item(51, 3).Value = '漫画/廃狼 原作/馬路まんじ キャラクター原案/かぼちゃ'
$ws.Cells.Item(51, 4).Value = '第8話（2）'
